$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4284.25
$ws.Range("I74").Value = 3600
$ws.Range("K74").Value = 3600
$ws.Range("M74").Value = -2664
$ws.Range("H77").Value = 4284.25
$ws.Range("I77").Value = 3600
$ws.Range("K77").Value = 18000
$ws.Range("M77").Value = -13320
$ws.Range("H98").Value = 508.2
$ws.Range("I98").Value = 231.33333
$ws.Range("K98").Value = 231.33333
$ws.Range("M98").Value = 1266.66667
$ws.Range("H100").Value = 1818.381
$ws.Range("I100").Value = 1641.5
$ws.Range("J100").Value = 2172.1428
$ws.Range("K100").Value = 1641.5
$ws.Range("L100").Value = 2172.1428
$ws.Range("M100").Value = -1100.5
$ws.Range("N100").Value = -3254.1428
$ws.Range("H122").Value = 508.2
$ws.Range("I122").Value = 231.33333
$ws.Range("K122").Value = 693.99999
$ws.Range("M122").Value = 1756.00001
$ws.Range("H135").Value = 1113.6428
$ws.Range("I135").Value = 659.46155
$ws.Range("K135").Value = 5935.15395
$ws.Range("M135").Value = -3400.15395
$ws.Range("H137").Value = 2528.4546
$ws.Range("I137").Value = 2034.5333
$ws.Range("J137").Value = 7467.6665
$ws.Range("K137").Value = 6103.5999
$ws.Range("L137").Value = 22402.9995
$ws.Range("M137").Value = -3553.5999
$ws.Range("N137").Value = -27502.9995
$ws.Range("H138").Value = 2100.7903
$ws.Range("I138").Value = 1564.8966
$ws.Range("J138").Value = 2571.7273
$ws.Range("K138").Value = 4694.6898
$ws.Range("L138").Value = 7715.1819
$ws.Range("M138").Value = 445.3101999999999
$ws.Range("N138").Value = -17995.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1833.2
$ws.Range("I2").Value = 1989
$ws.Range("J2").Value = 1599.5
$ws.Range("K2").Value = 1989
$ws.Range("L2").Value = 1599.5
$ws.Range("M2").Value = -1876
$ws.Range("N2").Value = -1825.5
$ws.Range("H32").Value = 452278.2
$ws.Range("I32").Value = 499721.44
$ws.Range("J32").Value = 25289
$ws.Range("K32").Value = 499721.44
$ws.Range("L32").Value = 25289
$ws.Range("M32").Value = -499434.44
$ws.Range("N32").Value = -25863
$ws.Range("H45").Value = 4073.2778
$ws.Range("I45").Value = 4033.2307
$ws.Range("J45").Value = 4177.4
$ws.Range("K45").Value = 4033.2307
$ws.Range("L45").Value = 4177.4
$ws.Range("M45").Value = -3656.2307
$ws.Range("N45").Value = -4931.4
$ws.Range("H61").Value = 2248.647
$ws.Range("I61").Value = 2115.68
$ws.Range("J61").Value = 2618
$ws.Range("K61").Value = 2115.68
$ws.Range("L61").Value = 2618
$ws.Range("M61").Value = -1903.68
$ws.Range("N61").Value = -3042
$ws.Range("H116").Value = 1833.2
$ws.Range("I116").Value = 1989
$ws.Range("J116").Value = 1599.5
$ws.Range("K116").Value = 1989
$ws.Range("L116").Value = 1599.5
$ws.Range("M116").Value = 305
$ws.Range("N116").Value = -6187.5
$ws.Range("H132").Value = 2759
$ws.Range("I132").Value = 1849.9736
$ws.Range("J132").Value = 4260.8696
$ws.Range("K132").Value = 5549.9208
$ws.Range("L132").Value = 12782.6088
$ws.Range("M132").Value = -3019.9208
$ws.Range("N132").Value = -17842.6088
$ws.Range("H136").Value = 2248.647
$ws.Range("I136").Value = 2115.68
$ws.Range("J136").Value = 2618
$ws.Range("K136").Value = 6347.039999999999
$ws.Range("L136").Value = 7854
$ws.Range("M136").Value = -3797.039999999999
$ws.Range("N136").Value = -12954
$ws.Range("H138").Value = 57944.855
$ws.Range("J138").Value = 57944.855
$ws.Range("L138").Value = 57944.855
$ws.Range("N138").Value = -68224.85500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1833.2
$ws.Range("I3").Value = 1989
$ws.Range("J3").Value = 1599.5
$ws.Range("K3").Value = 1989
$ws.Range("L3").Value = 1599.5
$ws.Range("M3").Value = -1875
$ws.Range("N3").Value = -1827.5
$ws.Range("H22").Value = 3493.5
$ws.Range("I22").Value = 3629.2727
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 3629.2727
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -3456.2727
$ws.Range("N22").Value = -2346
$ws.Range("H86").Value = 2069.2727
$ws.Range("I86").Value = 1959.0667
$ws.Range("J86").Value = 2305.4285
$ws.Range("K86").Value = 1959.0667
$ws.Range("L86").Value = 2305.4285
$ws.Range("M86").Value = -836.0667000000001
$ws.Range("N86").Value = -4551.4285
$ws.Range("H89").Value = 2069.2727
$ws.Range("I89").Value = 1959.0667
$ws.Range("J89").Value = 2305.4285
$ws.Range("K89").Value = 9795.333500000001
$ws.Range("L89").Value = 11527.1425
$ws.Range("M89").Value = -4179.333500000001
$ws.Range("N89").Value = -22759.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4326.526
$ws.Range("I31").Value = 1031.2
$ws.Range("J31").Value = 8820.151
$ws.Range("K31").Value = 1031.2
$ws.Range("L31").Value = 8820.151
$ws.Range("M31").Value = -736.2
$ws.Range("N31").Value = -9410.151
$ws.Range("H34").Value = 4326.526
$ws.Range("I34").Value = 1031.2
$ws.Range("J34").Value = 8820.151
$ws.Range("K34").Value = 1031.2
$ws.Range("L34").Value = 8820.151
$ws.Range("M34").Value = -829.2
$ws.Range("N34").Value = -9224.151
$ws.Range("H134").Value = 8737.4375
$ws.Range("I134").Value = 8342.429
$ws.Range("J134").Value = 11502.5
$ws.Range("K134").Value = 25027.287
$ws.Range("L134").Value = 34507.5
$ws.Range("M134").Value = -22492.287
$ws.Range("N134").Value = -39577.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1032.2858
$ws.Range("J131").Value = 1169.2941
$ws.Range("L131").Value = 3507.8823
$ws.Range("N131").Value = -13587.8823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 46312176
$ws.Range("J80").Value = 102249.75
$ws.Range("L80").Value = 102249.75
$ws.Range("N80").Value = -104245.75
$ws.Range("H83").Value = 46312176
$ws.Range("J83").Value = 102249.75
$ws.Range("L83").Value = 511248.75
$ws.Range("N83").Value = -521232.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 850
$ws.Range("J46").Value = 1150
$ws.Range("K46").Value = 850
$ws.Range("L46").Value = 1150
$ws.Range("M46").Value = -662
$ws.Range("N46").Value = -1526
$ws.Range("H132").Value = 2180.639
$ws.Range("I132").Value = 1569.619
$ws.Range("J132").Value = 3036.0667
$ws.Range("K132").Value = 4708.857
$ws.Range("L132").Value = 9108.2001
$ws.Range("M132").Value = -2178.857
$ws.Range("N132").Value = -14168.2001
$ws.Range("H136").Value = 11113001
$ws.Range("I136").Value = 2466.5
$ws.Range("J136").Value = 18520024
$ws.Range("K136").Value = 7399.5
$ws.Range("L136").Value = 55560072
$ws.Range("M136").Value = -4849.5
$ws.Range("N136").Value = -55565172

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5385.3335
$ws.Range("I81").Value = 6374
$ws.Range("J81").Value = 4149.5
$ws.Range("K81").Value = 12748
$ws.Range("L81").Value = 8299
$ws.Range("M81").Value = -11687
$ws.Range("N81").Value = -10421
$ws.Range("H84").Value = 5385.3335
$ws.Range("I84").Value = 6374
$ws.Range("J84").Value = 4149.5
$ws.Range("K84").Value = 63740
$ws.Range("L84").Value = 41495
$ws.Range("M84").Value = -58436
$ws.Range("N84").Value = -52103
$ws.Range("H136").Value = 1926.619
$ws.Range("I136").Value = 1417.1
$ws.Range("J136").Value = 3200.4167
$ws.Range("K136").Value = 4251.299999999999
$ws.Range("L136").Value = 9601.250100000001
$ws.Range("M136").Value = -1701.299999999999
$ws.Range("N136").Value = -14701.2501

